$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) columns with refreshed symbol-list data.
# Cells are forced to Text format ("@") so the values are stored as literal
# strings (matching the sheet's existing inline-string convention) rather
# than being auto-coerced into numbers/percentages by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "320.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.94%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.13%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.126"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.04%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08152"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.77%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.143"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.80%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.041"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.139"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.09%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9268"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.83%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1008"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.27%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1887"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.24%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09174"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.19%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03591"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.69%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09920"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.14%"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.94%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005693"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.08%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.451"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.29%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "15.86%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.57%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1310"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.45%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.058"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.01%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2188"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.98%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.90%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001243"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.70%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004729"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-7.15%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.18%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004500"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.40%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02031"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "11.25%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04993"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.68%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007847"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.02%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.28%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007807"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.00%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002095"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.91%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01215"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.17%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006468"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.36%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.15%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "19.50%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001900"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.13%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.15%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.15%"

